$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / Uncertainty: was a number (1.5), becomes the text "1.5".
# A leading apostrophe forces Excel to store it as text instead of
# re-parsing it back into a number; resetting the style afterwards keeps
# the cell on the workbook's default (unstyled) format.
$ws.Range("D2").Value = "'1.5"
$ws.Range("D2").Style = "Normal"

# Row 4 / Fraction of archaea: Value was text "0.2", becomes the number 0.2.
$ws.Range("B4").Value = 0.2

# Row 4 / Uncertainty: was text "2.3", becomes the number 2.2.
$ws.Range("D4").Value = 2.2

# Row 5 / Fraction of bacteria: Value was text "0.8", becomes the number 0.8.
$ws.Range("B5").Value = 0.8

# Row 5 / Uncertainty: was text "1.3", becomes the number 1.3.
$ws.Range("D5").Value = 1.3
